$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("39:41").Insert()
$ws.Rows("39:41").Clear()
